# Apply the "456a3b4" data refresh to 北京-漫展信息.xlsx
#
# Sheet "展览" (worksheet 1):
#   - SISP event (row 3) cancelled: title gets "（取消）" suffix, price cell
#     becomes the text "不可售" (not for sale).
#   - Various "want to go" counts (column F) bumped, a couple of "min price"
#     cells (column G) updated, and row 21 (Ming-Ri-Fang-Zhou only show)
#     also flips to "不可售".
#
# Sheet "演出" (worksheet 2):
#   - F4 bumped from 33 to 35.
#
# Sheet "全部类型" (worksheet 4):
#   - Mirrors the same edits as "展览" (it aggregates all categories), just
#     offset by two rows (row 5 instead of row 3, etc).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- 展览 ("Exhibitions") ----
$wsExhibit.Range("C3").Value = "北京·SISP动漫游戏嘉年华（取消）"
$wsExhibit.Range("G3").Value = "不可售"
$wsExhibit.Range("F4").Value = 1280
$wsExhibit.Range("F6").Value = 310
$wsExhibit.Range("F7").Value = 1113
$wsExhibit.Range("F9").Value = 6948
$wsExhibit.Range("G9").Value = 70.2
$wsExhibit.Range("F11").Value = 85
$wsExhibit.Range("F13").Value = 7843
$wsExhibit.Range("F16").Value = 5446
$wsExhibit.Range("F18").Value = 2322
$wsExhibit.Range("F19").Value = 979
$wsExhibit.Range("F21").Value = 274
$wsExhibit.Range("G21").Value = "不可售"
$wsExhibit.Range("F25").Value = 317
$wsExhibit.Range("F28").Value = 2096
$wsExhibit.Range("F32").Value = 48
$wsExhibit.Range("F33").Value = 547
$wsExhibit.Range("F36").Value = 1417
$wsExhibit.Range("F39").Value = 2153
$wsExhibit.Range("F40").Value = 2180

# ---- 演出 ("Shows") ----
$wsShow.Range("F4").Value = 35

# ---- 全部类型 ("All categories") ----
$wsAll.Range("C5").Value = "北京·SISP动漫游戏嘉年华（取消）"
$wsAll.Range("G5").Value = "不可售"
$wsAll.Range("F6").Value = 1280
$wsAll.Range("F9").Value = 310
$wsAll.Range("F10").Value = 1113
$wsAll.Range("F12").Value = 6948
$wsAll.Range("G12").Value = 70.2
$wsAll.Range("F14").Value = 85
$wsAll.Range("F16").Value = 7843
$wsAll.Range("F19").Value = 5446
$wsAll.Range("F21").Value = 2322
$wsAll.Range("F22").Value = 979
$wsAll.Range("F24").Value = 274
$wsAll.Range("G24").Value = "不可售"
$wsAll.Range("F29").Value = 35
$wsAll.Range("F30").Value = 317
$wsAll.Range("F33").Value = 2096
$wsAll.Range("F37").Value = 48
$wsAll.Range("F38").Value = 547
$wsAll.Range("F42").Value = 1417
$wsAll.Range("F45").Value = 2153
$wsAll.Range("F47").Value = 2180
